# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.757.98'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '2.102.17'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('D4').Value = "'1.008"
$ws.Range('E4').Value = '  +0.48%  '
$ws.Range('D5').Value = "'347.59"
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('D8').Value = "'0.4416"
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('D9').Value = "'53.91"
$ws.Range('E9').Value = '  +2.87%  '
$ws.Range('D10').Value = "'0.09388"
$ws.Range('E10').Value = '  +4.71%  '
$ws.Range('D11').Value = "'1.173"
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').Value = "'24.92"
$ws.Range('E12').Value = '  -2.83%  '
$ws.Range('D13').Value = '2.108.73'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').Value = "'6.824"
$ws.Range('E14').Value = '  +1.34%  '
$ws.Range('D15').Value = "'8.260"
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = "'102.76"
$ws.Range('E16').Value = '  +3.26%  '
$ws.Range('D17').Value = "'0.00001160"
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').Value = "'1.009"
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').Value = "'21.13"
$ws.Range('E19').Value = '  +1.61%  '
$ws.Range('D20').Value = "'0.06669"
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').Value = "'1.007"
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').Value = "'6.286"
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').Value = '29.799.86'
$ws.Range('E23').Value = '  -1.20%  '
$ws.Range('D24').Value = "'12.63"
$ws.Range('E24').Value = '  -1.00%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = "'2.316"
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.359.74'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('D28').Value = "'162.43"
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').Value = "'2.523"
$ws.Range('E29').Value = '  -0.64%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').Value = "'1.139"
$ws.Range('E31').Value = '  -3.10%  '
$ws.Range('D32').Value = "'1.733"
$ws.Range('E32').Value = '  +5.89%  '
$ws.Range('E33').Value = '  -1.12%  '
$ws.Range('D34').Value = "'6.221"
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').Value = "'3.949"
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('D36').Value = "'6.359"
$ws.Range('E36').Value = '  +6.77%  '
$ws.Range('D37').Value = "'10.49"
$ws.Range('E37').Value = '  +2.29%  '
$ws.Range('D38').Value = "'0.02587"
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('D39').Value = "'0.06754"
$ws.Range('E39').Value = '  -0.87%  '
$ws.Range('D40').Value = "'0.7000"
$ws.Range('E40').Value = '  +2.70%  '
$ws.Range('D41').Value = "'12.60"
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').Value = "'1.334"
$ws.Range('E42').Value = '  +4.02%  '
$ws.Range('D43').Value = "'0.2227"
$ws.Range('E43').Value = '  -3.06%  '
$ws.Range('D44').Value = "'0.6828"
$ws.Range('E44').Value = '  +6.88%  '
$ws.Range('D45').Value = "'14.43"
$ws.Range('E45').Value = '  +1.44%  '
$ws.Range('D46').Value = "'2.358"
$ws.Range('E46').Value = '  +2.70%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = "'3.636"
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = "'0.00000000356"
$ws.Range('E48').Value = '  -2.06%  '
$ws.Range('E49').Value = '  +4.50%  '
$ws.Range('D50').Value = "'1.222"
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').Value = "'81.63"
$ws.Range('E51').Value = '  -1.02%  '
